$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 9617713
$ws.Range("I80").Value = 17858510
$ws.Range("J80").Value = 3449.5
$ws.Range("K80").Value = 53575530
$ws.Range("L80").Value = 10348.5
$ws.Range("M80").Value = -53574532
$ws.Range("N80").Value = -12344.5
$ws.Range("H83").Value = 9617713
$ws.Range("I83").Value = 17858510
$ws.Range("J83").Value = 3449.5
$ws.Range("K83").Value = 160726590
$ws.Range("L83").Value = 31045.5
$ws.Range("M83").Value = -160721598
$ws.Range("N83").Value = -41029.5
$ws.Range("H111").Value = 542.2857
$ws.Range("I111").Value = 542.2857
$ws.Range("K111").Value = 1626.8571
$ws.Range("M111").Value = 1440.1429
$ws.Range("H112").Value = 1332.2963
$ws.Range("J112").Value = 1421.091
$ws.Range("L112").Value = 4263.272999999999
$ws.Range("N112").Value = -6479.272999999999
$ws.Range("H133").Value = 71779
$ws.Range("J133").Value = 71779
$ws.Range("L133").Value = 71779
$ws.Range("N133").Value = -81899
$ws.Range("H136").Value = 95000
$ws.Range("J136").Value = 95000
$ws.Range("L136").Value = 95000
$ws.Range("N136").Value = -105200
$ws.Range("H138").Value = 1612.663
$ws.Range("I138").Value = 856.8461
$ws.Range("J138").Value = 2674.8918
$ws.Range("K138").Value = 2570.5383
$ws.Range("L138").Value = 8024.6754
$ws.Range("M138").Value = 2569.4617
$ws.Range("N138").Value = -18304.6754
$ws.Range("H141").Value = 696.7045000000001
$ws.Range("I141").Value = 566.7
$ws.Range("K141").Value = 1700.1
$ws.Range("M141").Value = 3479.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 63390.5
$ws.Range("I102").Value = 63390.5
$ws.Range("K102").Value = 63390.5
$ws.Range("M102").Value = -61768.5
$ws.Range("H122").Value = 2787.926
$ws.Range("I122").Value = 2310.3
$ws.Range("K122").Value = 6930.900000000001
$ws.Range("M122").Value = -4480.900000000001
$ws.Range("H132").Value = 403046.75
$ws.Range("I132").Value = 424394.8
$ws.Range("J132").Value = 11665.667
$ws.Range("K132").Value = 1273184.4
$ws.Range("L132").Value = 34997.001
$ws.Range("M132").Value = -1270654.4
$ws.Range("N132").Value = -40057.001
$ws.Range("H134").Value = 88754.28999999999
$ws.Range("J134").Value = 88723.336
$ws.Range("L134").Value = 88723.336
$ws.Range("N134").Value = -98863.336
$ws.Range("H139").Value = 79950.836
$ws.Range("J139").Value = 79950.836
$ws.Range("L139").Value = 79950.836
$ws.Range("N139").Value = -90230.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 24640.777
$ws.Range("I99").Value = 24251.738
$ws.Range("J99").Value = 26877.75
$ws.Range("K99").Value = 24251.738
$ws.Range("L99").Value = 26877.75
$ws.Range("M99").Value = -22753.738
$ws.Range("N99").Value = -29873.75
$ws.Range("H107").Value = 3196.1177
$ws.Range("I107").Value = 3086
$ws.Range("K107").Value = 3086
$ws.Range("M107").Value = -1166
$ws.Range("H140").Value = 99219.39999999999
$ws.Range("J140").Value = 99219.39999999999
$ws.Range("L140").Value = 99219.39999999999
$ws.Range("N140").Value = -109579.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 104466.31
$ws.Range("I31").Value = 159759
$ws.Range("J31").Value = 30742.732
$ws.Range("K31").Value = 159759
$ws.Range("L31").Value = 30742.732
$ws.Range("M31").Value = -159464
$ws.Range("N31").Value = -31332.732
$ws.Range("H34").Value = 104466.31
$ws.Range("I34").Value = 159759
$ws.Range("J34").Value = 30742.732
$ws.Range("K34").Value = 159759
$ws.Range("L34").Value = 30742.732
$ws.Range("M34").Value = -159557
$ws.Range("N34").Value = -31146.732
$ws.Range("H68").Value = 87695
$ws.Range("J68").Value = 87695
$ws.Range("L68").Value = 87695
$ws.Range("N68").Value = -89193
$ws.Range("H71").Value = 87695
$ws.Range("J71").Value = 87695
$ws.Range("L71").Value = 263085
$ws.Range("N71").Value = -270573
$ws.Range("H99").Value = 3768.4243
$ws.Range("I99").Value = 3654.7778
$ws.Range("J99").Value = 3904.8
$ws.Range("K99").Value = 3654.7778
$ws.Range("L99").Value = 3904.8
$ws.Range("M99").Value = -2156.7778
$ws.Range("N99").Value = -6900.8
$ws.Range("H126").Value = 3768.4243
$ws.Range("I126").Value = 3654.7778
$ws.Range("J126").Value = 3904.8
$ws.Range("K126").Value = 10964.3334
$ws.Range("L126").Value = 11714.4
$ws.Range("M126").Value = -8494.3334
$ws.Range("N126").Value = -16654.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2969.32
$ws.Range("I102").Value = 1860.0625
$ws.Range("J102").Value = 4941.3335
$ws.Range("K102").Value = 1860.0625
$ws.Range("L102").Value = 4941.3335
$ws.Range("M102").Value = -238.0625
$ws.Range("N102").Value = -8185.3335
$ws.Range("H107").Value = 22966.875
$ws.Range("I107").Value = 25504.785
$ws.Range("J107").Value = 5201.5
$ws.Range("K107").Value = 25504.785
$ws.Range("L107").Value = 5201.5
$ws.Range("M107").Value = -23584.785
$ws.Range("N107").Value = -9041.5
$ws.Range("H122").Value = 9908.888999999999
$ws.Range("I122").Value = 4036.4
$ws.Range("J122").Value = 17249.5
$ws.Range("K122").Value = 12109.2
$ws.Range("L122").Value = 51748.5
$ws.Range("M122").Value = -9659.200000000001
$ws.Range("N122").Value = -56648.5
$ws.Range("H132").Value = 275040
$ws.Range("I132").Value = 345373.78
$ws.Range("J132").Value = 1519.7778
$ws.Range("K132").Value = 1036121.34
$ws.Range("L132").Value = 4559.3334
$ws.Range("M132").Value = -1033591.34
$ws.Range("N132").Value = -9619.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3804.7727
$ws.Range("I40").Value = 3512.25
$ws.Range("J40").Value = 4584.8335
$ws.Range("K40").Value = 3512.25
$ws.Range("L40").Value = 4584.8335
$ws.Range("M40").Value = -3376.25
$ws.Range("N40").Value = -4856.8335
$ws.Range("H55").Value = 1034.2084
$ws.Range("I55").Value = 222.44444
$ws.Range("J55").Value = 1521.2667
$ws.Range("K55").Value = 222.44444
$ws.Range("L55").Value = 1521.2667
$ws.Range("M55").Value = -49.44443999999999
$ws.Range("N55").Value = -1867.2667
$ws.Range("H93").Value = 2496.4736
$ws.Range("I93").Value = 2445.5
$ws.Range("J93").Value = 2639.2
$ws.Range("K93").Value = 2445.5
$ws.Range("L93").Value = 2639.2
$ws.Range("M93").Value = -1197.5
$ws.Range("N93").Value = -5135.2
$ws.Range("H100").Value = 11818.363
$ws.Range("I100").Value = 2799.8
$ws.Range("J100").Value = 19333.834
$ws.Range("K100").Value = 2799.8
$ws.Range("L100").Value = 19333.834
$ws.Range("M100").Value = -2258.8
$ws.Range("N100").Value = -20415.834
$ws.Range("H122").Value = 2759.7778
$ws.Range("I122").Value = 2462.5518
$ws.Range("J122").Value = 3991.1428
$ws.Range("K122").Value = 7387.655400000001
$ws.Range("L122").Value = 11973.4284
$ws.Range("M122").Value = -4937.655400000001
$ws.Range("N122").Value = -16873.4284
$ws.Range("H132").Value = 737542.3
$ws.Range("I132").Value = 753532.4
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2260597.2
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2258067.2
$ws.Range("N132").Value = -11060
$ws.Range("H136").Value = 23866.361
$ws.Range("I136").Value = 2982.392
$ws.Range("J136").Value = 130374.6
$ws.Range("K136").Value = 8947.175999999999
$ws.Range("L136").Value = 391123.8
$ws.Range("M136").Value = -6397.175999999999
$ws.Range("N136").Value = -396223.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2302.075
$ws.Range("I107").Value = 1356.7667
$ws.Range("J107").Value = 5138
$ws.Range("K107").Value = 4070.300099999999
$ws.Range("L107").Value = 15414
$ws.Range("M107").Value = -2150.300099999999
$ws.Range("N107").Value = -19254
$ws.Range("H122").Value = 1594.1746
$ws.Range("I122").Value = 1536.614
$ws.Range("J122").Value = 2141
$ws.Range("K122").Value = 4609.842000000001
$ws.Range("L122").Value = 6423
$ws.Range("M122").Value = -2159.842000000001
$ws.Range("N122").Value = -11323
